# Update the "Assigned" (D) and "Drivers" (E) columns of the schedule sheet
# to reflect the new roster assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "Paul, Thor Waguespack, Alejandro E. Ulvert";            E = "Paul, Thor Waguespack, Alejandro E. Ulvert" },
    @{ Row = 3;  D = "Ben Kairouz, Alexander, Jamari Pitchford";              E = "Ben Kairouz, Alexander, Jamari Pitchford" },
    @{ Row = 4;  D = "Jaxon, Alejandro Espinosa, George Ryckman";             E = "Jaxon, Alejandro Espinosa, George Ryckman" },
    @{ Row = 5;  D = "Jack Mogelof, Jake Dieterich, Adi";                     E = "Jack Mogelof, Adi" },
    @{ Row = 6;  D = "Gabe Heller, Alejandro L";                              E = "Gabe Heller" },
    @{ Row = 7;  D = "Henry";                                                 E = "" },
    @{ Row = 8;  D = "Noah Yaffe";                                            E = "Noah Yaffe" },
    @{ Row = 9;  D = "Jack Mogelof, Jake Dieterich, Adi";                     E = "Jack Mogelof, Adi" },
    @{ Row = 10; D = "Ben Kairouz, Paul, Alejandro Espinosa";                 E = "Ben Kairouz, Paul, Alejandro Espinosa" },
    @{ Row = 11; D = "Ezana, Edu, Blake Steel";                               E = "Edu, Blake Steel" },
    @{ Row = 12; D = "Alexander, Jamari Pitchford, George Ryckman";           E = "Alexander, Jamari Pitchford, George Ryckman" },
    @{ Row = 14; D = "Henry";                                                 E = "" },
    @{ Row = 15; D = "";                                                      E = "" },
    @{ Row = 16; D = "Matheo, Alejandro L, Noah Yaffe";                       E = "Matheo, Noah Yaffe" },
    @{ Row = 17; D = "Gabe Heller, Ben Kairouz, Thor Waguespack";             E = "Gabe Heller, Ben Kairouz, Thor Waguespack" },
    @{ Row = 18; D = "Jaxon, Kamsi, Josh Greene";                             E = "Jaxon, Josh Greene" },
    @{ Row = 19; D = "Edu, Ali Awada, Harry Corbin";                          E = "Edu, Harry Corbin" },
    @{ Row = 20; D = "Ezana, Blake Steel";                                    E = "Blake Steel" },
    @{ Row = 22; D = "Henry";                                                 E = "" },
    @{ Row = 23; D = "Kamsi, Josh Greene, Matheo";                            E = "Josh Greene, Matheo" },
    @{ Row = 24; D = "Paul, Jamari Pitchford, Thor Waguespack";               E = "Paul, Jamari Pitchford, Thor Waguespack" },
    @{ Row = 25; D = "Gabe Heller, Ezana, Blake Steel";                       E = "Gabe Heller, Blake Steel" },
    @{ Row = 26; D = "Jack Mogelof, George Ryckman, Alejandro E. Ulvert";     E = "Jack Mogelof, George Ryckman, Alejandro E. Ulvert" },
    @{ Row = 28; D = "Henry";                                                 E = "" },
    @{ Row = 29; D = "";                                                      E = "" },
    @{ Row = 30; D = "Jaxon, Jake Dieterich, Alejandro E. Ulvert";            E = "Jaxon, Alejandro E. Ulvert" },
    @{ Row = 31; D = "Jamari Pitchford, Thor Waguespack, Ali Awada";          E = "Jamari Pitchford, Thor Waguespack" },
    @{ Row = 32; D = "Ben Kairouz, George Ryckman, Adi";                      E = "Ben Kairouz, George Ryckman, Adi" },
    @{ Row = 33; D = "Jack Mogelof, Alejandro Espinosa, Harry Corbin";        E = "Jack Mogelof, Alejandro Espinosa, Harry Corbin" },
    @{ Row = 35; D = "Kamsi";                                                 E = "" },
    @{ Row = 36; D = "";                                                      E = "" },
    @{ Row = 37; D = "Josh Greene, Adi, Harry Corbin";                        E = "Josh Greene, Adi, Harry Corbin" },
    @{ Row = 38; D = "Jaxon, Paul, Alejandro L";                              E = "Jaxon, Paul" },
    @{ Row = 39; D = "Alexander, Edu, Jake Dieterich";                        E = "Alexander, Edu" },
    @{ Row = 40; D = "Alejandro Espinosa, Ali Awada, Alejandro E. Ulvert";    E = "Alejandro Espinosa, Alejandro E. Ulvert" },
    @{ Row = 43; D = "Matheo, Noah Yaffe";                                    E = "Matheo, Noah Yaffe" },
    @{ Row = 44; D = "Josh Greene, Alejandro L, Harry Corbin";                E = "Josh Greene, Harry Corbin" },
    @{ Row = 45; D = "Ezana, Alexander, Edu";                                 E = "Alexander, Edu" },
    @{ Row = 46; D = "Gabe Heller, Ali Awada, Noah Yaffe";                    E = "Gabe Heller, Noah Yaffe" },
    @{ Row = 47; D = "Kamsi, Blake Steel";                                    E = "Blake Steel" },
    @{ Row = 48; D = "Matheo";                                                E = "Matheo" },
    @{ Row = 50; D = "";                                                      E = "" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

Write-Host "Updated $($updates.Count) rows"
